$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 31 with the latest bulletin data
$ws.Range("A31").Value = 1109
$ws.Range("B31").Value = 462
$ws.Range("C31").Value = 94
$ws.Range("D31").Value = 145
$ws.Range("E31").Value = 122
$ws.Range("F31").Formula = "=+A31-SUM(B31:E31)"

# Select column G (whole column) as the active selection, matching the saved view state
$ws.Columns.Item(7).Select()
